$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.807785987854004
$ws.Range("B1").Value = 4.002347469329834
$ws.Range("C1").Value = 3.035057783126831
$ws.Range("D1").Value = 2.659415483474731
$ws.Range("E1").Value = 1.868397831916809
